$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.529.23'
$ws.Range("E2").Value = '  +2.85%  '

$ws.Range("D3").Value = '2.353.28'
$ws.Range("E3").Value = '  +6.09%  '

$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").Value = '314.82'
$ws.Range("E5").Value = '  +6.09%  '

$ws.Range("D6").Value = '110.35'
$ws.Range("E6").Value = '  +2.46%  '

$ws.Range("E7").Value = '  +3.78%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("E9").Value = '  +6.60%  '

$ws.Range("D10").Value = '43.32'
$ws.Range("E10").Value = '  -0.66%  '

$ws.Range("E11").Value = '  +3.07%  '

$ws.Range("D12").Value = '8.85'
$ws.Range("E12").Value = '  +0.94%  '

$ws.Range("D13").Value = '1.05'
$ws.Range("E13").Value = '  +6.85%  '

$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +2.47%  '

$ws.Range("D15").Value = '16.42'
$ws.Range("E15").Value = '  +9.86%  '

$ws.Range("D16").Value = '2.709.09'
$ws.Range("E16").Value = '  +6.20%  '

$ws.Range("D17").Value = '2.423.60'
$ws.Range("E17").Value = '  +8.69%  '

$ws.Range("D18").Value = '43.486.39'
$ws.Range("E18").Value = '  +2.99%  '

$ws.Range("E19").Value = '  +3.85%  '

$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  -1.73%  '

$ws.Range("D21").Value = '75.66'
$ws.Range("E21").Value = '  +4.55%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("E23").Value = '  +12.65%  '

$ws.Range("D24").Value = '255.60'
$ws.Range("E24").Value = '  +12.05%  '

$ws.Range("D25").Value = '9.14'
$ws.Range("E25").Value = '  +0.86%  '

$ws.Range("D26").Value = '12.07'
$ws.Range("E26").Value = '  +4.29%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").Value = '39.28'
$ws.Range("E28").Value = '  +2.87%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +1.03%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '22.42'
$ws.Range("E30").Value = '  +6.98%  '

$ws.Range("D31").Value = '173.92'
$ws.Range("E31").Value = '  +0.10%  '

$ws.Range("D32").Value = '3.18'
$ws.Range("E32").Value = '  -0.82%  '

$ws.Range("D33").Value = '0.0933'
$ws.Range("E33").Value = '  +4.16%  '

$ws.Range("E34").Value = '  +7.75%  '

$ws.Range("D35").Value = '0.133'
$ws.Range("E35").Value = '  +5.64%  '

$ws.Range("E36").Value = '  -1.94%  '

$ws.Range("D37").Value = '4.19'
$ws.Range("E37").Value = '  -3.54%  '

$ws.Range("E38").Value = '  +2.15%  '

$ws.Range("D39").Value = '0.105'
$ws.Range("E39").Value = '  +2.30%  '

$ws.Range("D40").Value = '2.71'
$ws.Range("E40").Value = '  +12.39%  '

$ws.Range("D41").Value = '72.73'
$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("E42").Value = '  +14.58%  '

$ws.Range("D43").Value = '0.235'
$ws.Range("E43").Value = '  +1.53%  '

$ws.Range("D44").Value = '12.86'
$ws.Range("E44").Value = '  +2.50%  '

$ws.Range("E45").Value = '  +0.20%  '

$ws.Range("D46").Value = '5.65'
$ws.Range("E46").Value = '  +4.68%  '

$ws.Range("D47").Value = '9.31'
$ws.Range("E47").Value = '  +10.90%  '

$ws.Range("D48").Value = '111.53'
$ws.Range("E48").Value = '  +8.13%  '

$ws.Range("D49").Value = '1.32'
$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  +2.91%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '70.48'
$ws.Range("E51").Value = '  +5.43%  '
